# Remove the trailing "Ver no Jupiter..." / copyright boilerplate block
# (and the blank paragraph that preceded it) that followed the
# "Bibliografia" section content, while keeping the blank paragraph and
# page-break paragraph that close out the document.

$d = $word.ActiveDocument

# Locate the paragraph whose text is the blank line right after the
# "Artigos sobre metodologias ativas..." bibliography sentence, and the
# paragraph containing the trailing copyright notice, by scanning the
# paragraph collection for their known text content. This avoids relying
# on hard-coded paragraph indices.

$blankIndex = -1
$jupiterIndex = -1
$copyrightIndex = -1

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "Ver no Jupiter") {
        $jupiterIndex = $i
        $blankIndex = $i - 1
    }
    if ($t -match "Powered by Jekyll") {
        $copyrightIndex = $i
    }
}

if ($jupiterIndex -gt 0 -and $copyrightIndex -gt 0 -and $blankIndex -gt 0) {
    $startRange = $d.Paragraphs.Item($blankIndex).Range.Start
    $endRange = $d.Paragraphs.Item($copyrightIndex).Range.End
    $r = $d.Range($startRange, $endRange)
    $r.Delete()
}
